$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.000.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.265.81"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.43"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.656"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.59"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.08%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.15"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.66"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.602.26"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.63"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.15"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.03%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.263.91"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.864.48"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0984"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.01"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.35"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.48"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.36"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +20.99%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.23"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.70%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.55%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "174.02"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.35%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.12"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.02%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.96"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0687"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.96"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.72"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.98%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0255"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.49%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.76"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000222"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.89"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.27"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.42"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0953"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.20"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.456.69"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.33"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.00"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -6.10%  "
